$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in lat/lon values for rows 4-9 (columns A and B)
$ws.Range("A4").Value = -8.3444
$ws.Range("B4").Value = 114.6272

$ws.Range("A5").Value = -8.722749
$ws.Range("B5").Value = 115.1697

$ws.Range("A6").Value = -8.746536
$ws.Range("B6").Value = 115.211375

$ws.Range("A7").Value = -8.46083333
$ws.Range("B7").Value = 115.1396667

$ws.Range("A8").Value = -8.41598
$ws.Range("B8").Value = 115.4201

$ws.Range("A9").Value = -8.54216
$ws.Range("B9").Value = 115.441

# Update the active cell selection to B9
$ws.Range("B9").Select()
